$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Completed Date" column (E) was missing values for rows 26-28.
# Fill them in with 45117 (2023-07-10), reusing the date style already
# used throughout the rest of column E (copy format from E25).
$ws.Range("E26").Value = 45117
$ws.Range("E27").Value = 45117
$ws.Range("E28").Value = 45117

$ws.Range("E25").Copy()
$ws.Range("E26:E28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to the last-edited cell, E28.
$ws.Range("E28").Select()
